# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the content is
# now "Ready for handoff" (instead of "Handed back: in sync with en-US"),
# and refreshes the related handoff timestamps.

$wb = $excel.ActiveWorkbook

$Overview = $wb.Worksheets.Item(1)   # "Overview" sheet
$ZhCn     = $wb.Worksheets.Item(2)   # "zh-cn" sheet
$DeDe     = $wb.Worksheets.Item(3)   # "de-de" sheet

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$Overview.Range("E2").Value = "Ready for handoff"
$Overview.Range("F2").Value = "Ready for handoff"
$ZhCn.Range("C2").Value = "Ready for handoff"
$DeDe.Range("C2").Value = "Ready for handoff"

# --- Timestamps ---
# Overview "Latest HO Xliff Generate Date" and de-de "Latest Handoff Datetime"
$Overview.Range("G2").Value = "2016-08-25 04:57:31"
$DeDe.Range("H2").Value = "2016-08-25 04:57:31"

# zh-cn "Latest Handoff Datetime"
$ZhCn.Range("H2").Value = "2016-08-25 04:57:27"

# --- Column widths narrowed for the Status / language-status columns ---
$Overview.Columns.Item(5).ColumnWidth = 16.25   # column E (zh-cn)
$Overview.Columns.Item(6).ColumnWidth = 16.25   # column F (de-de)
$ZhCn.Columns.Item(3).ColumnWidth = 16.25        # column C (Status)
$DeDe.Columns.Item(3).ColumnWidth = 16.25        # column C (Status)
